# Updates the cached "datetimeFigureOut" date text on the slide master and
# every slide layout (19-2-2021 -> 19-6-2022), and adds two small "+"/"-"
# textboxes to slide 1, matching the target edit.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached date placeholder text across the master + layouts
# ---------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "19-2-2021") {
                $sh.TextFrame.TextRange.Text = "19-6-2022"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DateShapes $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Add the two new "+" / "-" textboxes on slide 1
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

# EMU -> point conversion (the shape-position API works in points = EMU/12700)
$plusLeft   = 5557421 / 12700
$plusTop    = 2175029 / 12700
$plusWidth  = 300082 / 12700
$plusHeight = 369332 / 12700

$plusBox = $s.Shapes.AddTextbox(1, $plusLeft, $plusTop, $plusWidth, $plusHeight)
$plusBox.Name = "TextBox 1"
$plusBox.TextFrame.WordWrap = $false
$plusBox.TextFrame.AutoSize = 1
$plusBox.Fill.Visible = $false
$plusRange = $plusBox.TextFrame.TextRange
$plusRange.Text = "+"
$plusRange.Font.Color.ObjectThemeColor = 10

$minusLeft   = 5767851 / 12700
$minusTop    = 2175029 / 12700
$minusWidth  = 255198 / 12700
$minusHeight = 369332 / 12700

$minusBox = $s.Shapes.AddTextbox(1, $minusLeft, $minusTop, $minusWidth, $minusHeight)
$minusBox.Name = "TextBox 6"
$minusBox.TextFrame.WordWrap = $false
$minusBox.TextFrame.AutoSize = 1
$minusBox.Fill.Visible = $false
$minusRange = $minusBox.TextFrame.TextRange
$minusRange.Text = "-"
$minusRange.Font.Bold = $true
$minusRange.Font.Color.RGB = 255
